$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = [double]"2"
$ws.Range("G2").Value = [double]"47.0283035"
$ws.Range("H2").Value = [double]"94.056607"
$ws.Range("I2").Value = [double]"0.1170896029811303"
$ws.Range("J2").Value = [double]"0.08213676148506427"
$ws.Range("K2").Value = [double]"2"
$ws.Range("M2").Value = [double]"35.789624"
$ws.Range("N2").Value = [double]"71.57924800000001"
$ws.Range("O2").Value = [double]"0.258139457682779"
$ws.Range("P2").Value = [double]"0.1993778771086309"
$ws.Range("Q2").Value = [double]"1683.125299622884"
$ws.Range("R2").Value = [double]"6732.501198491536"
$ws.Range("S2").Value = [double]"0.03022544661384088"
$ws.Range("T2").Value = [double]"0.01637625313747007"

# Row 3
$ws.Range("E3").Value = [double]"2"
$ws.Range("G3").Value = [double]"47.0283035"
$ws.Range("H3").Value = [double]"94.056607"
$ws.Range("I3").Value = [double]"0.1170896029811303"
$ws.Range("J3").Value = [double]"0.08213676148506427"
$ws.Range("K3").Value = [double]"3"
$ws.Range("M3").Value = [double]"60.113367"
$ws.Range("N3").Value = [double]"180.340101"
$ws.Range("O3").Value = [double]"0.4335790718803266"
$ws.Range("P3").Value = [double]"0.5023219368682956"
$ws.Range("Q3").Value = [double]"2827.029667682885"
$ws.Range("R3").Value = [double]"16962.17800609731"
$ws.Range("S3").Value = [double]"0.05076760138739442"
$ws.Range("T3").Value = [double]"0.04125909711726671"

# Row 4
$ws.Range("E4").Value = [double]"2"
$ws.Range("G4").Value = [double]"47.0283035"
$ws.Range("H4").Value = [double]"94.056607"
$ws.Range("I4").Value = [double]"0.1170896029811303"
$ws.Range("J4").Value = [double]"0.08213676148506427"
$ws.Range("K4").Value = [double]"3"
$ws.Range("M4").Value = [double]"0.2072186666666667"
$ws.Range("N4").Value = [double]"0.621656"
$ws.Range("O4").Value = [double]"0.001494603973349423"
$ws.Range("P4").Value = [double]"0.001731569652308208"
$ws.Range("Q4").Value = [double]"9.745142346865332"
$ws.Range("R4").Value = [double]"58.470854081192"
$ws.Range("S4").Value = [double]"0.0001750025858535038"
$ws.Range("T4").Value = [double]"0.000142225523526415"

# Row 5
$ws.Range("E5").Value = [double]"2"
$ws.Range("G5").Value = [double]"47.0283035"
$ws.Range("H5").Value = [double]"94.056607"
$ws.Range("I5").Value = [double]"0.1170896029811303"
$ws.Range("J5").Value = [double]"0.08213676148506427"
$ws.Range("K5").Value = [double]"3"
$ws.Range("M5").Value = [double]"21.40334366666667"
$ws.Range("N5").Value = [double]"64.210031"
$ws.Range("O5").Value = [double]"0.1543756795743782"
$ws.Range("P5").Value = [double]"0.178851553034748"
$ws.Range("Q5").Value = [double]"1006.562941870803"
$ws.Range("R5").Value = [double]"6039.377651224817"
$ws.Range("S5").Value = [double]"0.01807578703130613"
$ws.Range("T5").Value = [double]"0.01469028735284842"

# Row 6
$ws.Range("E6").Value = [double]"2"
$ws.Range("G6").Value = [double]"47.0283035"
$ws.Range("H6").Value = [double]"94.056607"
$ws.Range("I6").Value = [double]"0.1170896029811303"
$ws.Range("J6").Value = [double]"0.08213676148506427"
$ws.Range("K6").Value = [double]"2"
$ws.Range("M6").Value = [double]"21.1309775"
$ws.Range("N6").Value = [double]"42.261955"
$ws.Range("O6").Value = [double]"0.1524111868891667"
$ws.Range("P6").Value = [double]"0.1177170633360173"
$ws.Range("Q6").Value = [double]"993.7540231216713"
$ws.Range("R6").Value = [double]"3975.016092486685"
$ws.Range("S6").Value = [double]"0.01784576536273538"
$ws.Range("T6").Value = [double]"0.009668898353952656"

# Row 7
$ws.Range("E7").Value = [double]"3"
$ws.Range("G7").Value = [double]"12.33126566666667"
$ws.Range("H7").Value = [double]"36.993797"
$ws.Range("I7").Value = [double]"0.03070200057641551"
$ws.Range("J7").Value = [double]"0.03230555276798244"
$ws.Range("K7").Value = [double]"2"
$ws.Range("M7").Value = [double]"35.789624"
$ws.Range("N7").Value = [double]"71.57924800000001"
$ws.Range("O7").Value = [double]"0.258139457682779"
$ws.Range("P7").Value = [double]"0.1993778771086309"
$ws.Range("Q7").Value = [double]"441.3313616541094"
$ws.Range("R7").Value = [double]"2647.988169924656"
$ws.Range("S7").Value = [double]"0.007925397778572266"
$ws.Range("T7").Value = [double]"0.006441012529701193"

# Row 8
$ws.Range("E8").Value = [double]"3"
$ws.Range("G8").Value = [double]"12.33126566666667"
$ws.Range("H8").Value = [double]"36.993797"
$ws.Range("I8").Value = [double]"0.03070200057641551"
$ws.Range("J8").Value = [double]"0.03230555276798244"
$ws.Range("K8").Value = [double]"3"
$ws.Range("M8").Value = [double]"60.113367"
$ws.Range("N8").Value = [double]"180.340101"
$ws.Range("O8").Value = [double]"0.4335790718803266"
$ws.Range("P8").Value = [double]"0.5023219368682956"
$ws.Range("Q8").Value = [double]"741.2738985948331"
$ws.Range("R8").Value = [double]"6671.465087353497"
$ws.Range("S8").Value = [double]"0.01331174491479149"
$ws.Range("T8").Value = [double]"0.01622778783801387"

# Row 9
$ws.Range("E9").Value = [double]"3"
$ws.Range("G9").Value = [double]"12.33126566666667"
$ws.Range("H9").Value = [double]"36.993797"
$ws.Range("I9").Value = [double]"0.03070200057641551"
$ws.Range("J9").Value = [double]"0.03230555276798244"
$ws.Range("K9").Value = [double]"3"
$ws.Range("M9").Value = [double]"0.2072186666666667"
$ws.Range("N9").Value = [double]"0.621656"
$ws.Range("O9").Value = [double]"0.001494603973349423"
$ws.Range("P9").Value = [double]"0.001731569652308208"
$ws.Range("Q9").Value = [double]"2.555268429759111"
$ws.Range("R9").Value = [double]"22.997415867832"
$ws.Range("S9").Value = [double]"4.588733205128689E-05"
$ws.Range("T9").Value = [double]"5.593931477407983E-05"

# Row 10
$ws.Range("E10").Value = [double]"3"
$ws.Range("G10").Value = [double]"12.33126566666667"
$ws.Range("H10").Value = [double]"36.993797"
$ws.Range("I10").Value = [double]"0.03070200057641551"
$ws.Range("J10").Value = [double]"0.03230555276798244"
$ws.Range("K10").Value = [double]"3"
$ws.Range("M10").Value = [double]"21.40334366666667"
$ws.Range("N10").Value = [double]"64.210031"
$ws.Range("O10").Value = [double]"0.1543756795743782"
$ws.Range("P10").Value = [double]"0.178851553034748"
$ws.Range("Q10").Value = [double]"263.9303169086342"
$ws.Range("R10").Value = [double]"2375.372852177707"
$ws.Range("S10").Value = [double]"0.004739642203277095"
$ws.Range("T10").Value = [double]"0.00577789828419966"

# Row 11
$ws.Range("E11").Value = [double]"3"
$ws.Range("G11").Value = [double]"12.33126566666667"
$ws.Range("H11").Value = [double]"36.993797"
$ws.Range("I11").Value = [double]"0.03070200057641551"
$ws.Range("J11").Value = [double]"0.03230555276798244"
$ws.Range("K11").Value = [double]"2"
$ws.Range("M11").Value = [double]"21.1309775"
$ws.Range("N11").Value = [double]"42.261955"
$ws.Range("O11").Value = [double]"0.1524111868891667"
$ws.Range("P11").Value = [double]"0.1177170633360173"
$ws.Range("Q11").Value = [double]"260.5716973488558"
$ws.Range("R11").Value = [double]"1563.430184093135"
$ws.Range("S11").Value = [double]"0.004679328347723367"
$ws.Range("T11").Value = [double]"0.003802914801293637"

# Row 12
$ws.Range("E12").Value = [double]"3"
$ws.Range("G12").Value = [double]"129.0494106666667"
$ws.Range("H12").Value = [double]"387.148232"
$ws.Range("I12").Value = [double]"0.3213031969122349"
$ws.Range("J12").Value = [double]"0.3380847237148192"
$ws.Range("K12").Value = [double]"2"
$ws.Range("M12").Value = [double]"35.789624"
$ws.Range("N12").Value = [double]"71.57924800000001"
$ws.Range("O12").Value = [double]"0.258139457682779"
$ws.Range("P12").Value = [double]"0.1993778771086309"
$ws.Range("Q12").Value = [double]"4618.629885181589"
$ws.Range("R12").Value = [double]"27711.77931108954"
$ws.Range("S12").Value = [double]"0.08294103300266745"
$ws.Range("T12").Value = [double]"0.06740661449711864"

# Row 13
$ws.Range("E13").Value = [double]"3"
$ws.Range("G13").Value = [double]"129.0494106666667"
$ws.Range("H13").Value = [double]"387.148232"
$ws.Range("I13").Value = [double]"0.3213031969122349"
$ws.Range("J13").Value = [double]"0.3380847237148192"
$ws.Range("K13").Value = [double]"3"
$ws.Range("M13").Value = [double]"60.113367"
$ws.Range("N13").Value = [double]"180.340101"
$ws.Range("O13").Value = [double]"0.4335790718803266"
$ws.Range("P13").Value = [double]"0.5023219368682956"
$ws.Range("Q13").Value = [double]"7757.594584539048"
$ws.Range("R13").Value = [double]"69818.35126085143"
$ws.Range("S13").Value = [double]"0.1393103419093886"
$ws.Range("T13").Value = [double]"0.1698273732420106"

# Row 14
$ws.Range("E14").Value = [double]"3"
$ws.Range("G14").Value = [double]"129.0494106666667"
$ws.Range("H14").Value = [double]"387.148232"
$ws.Range("I14").Value = [double]"0.3213031969122349"
$ws.Range("J14").Value = [double]"0.3380847237148192"
$ws.Range("K14").Value = [double]"3"
$ws.Range("M14").Value = [double]"0.2072186666666667"
$ws.Range("N14").Value = [double]"0.621656"
$ws.Range("O14").Value = [double]"0.001494603973349423"
$ws.Range("P14").Value = [double]"0.001731569652308208"
$ws.Range("Q14").Value = [double]"26.74144681246577"
$ws.Range("R14").Value = [double]"240.673021312192"
$ws.Range("S14").Value = [double]"0.0004802210347548983"
$ws.Range("T14").Value = [double]"0.0005854172474935861"

# Row 15
$ws.Range("E15").Value = [double]"3"
$ws.Range("G15").Value = [double]"129.0494106666667"
$ws.Range("H15").Value = [double]"387.148232"
$ws.Range("I15").Value = [double]"0.3213031969122349"
$ws.Range("J15").Value = [double]"0.3380847237148192"
$ws.Range("K15").Value = [double]"3"
$ws.Range("M15").Value = [double]"21.40334366666667"
$ws.Range("N15").Value = [double]"64.210031"
$ws.Range("O15").Value = [double]"0.1543756795743782"
$ws.Range("P15").Value = [double]"0.178851553034748"
$ws.Range("Q15").Value = [double]"2762.088886479466"
$ws.Range("R15").Value = [double]"24858.79997831519"
$ws.Range("S15").Value = [double]"0.04960139937274651"
$ws.Range("T15").Value = [double]"0.0604669778937191"

# Row 16
$ws.Range("E16").Value = [double]"3"
$ws.Range("G16").Value = [double]"129.0494106666667"
$ws.Range("H16").Value = [double]"387.148232"
$ws.Range("I16").Value = [double]"0.3213031969122349"
$ws.Range("J16").Value = [double]"0.3380847237148192"
$ws.Range("K16").Value = [double]"2"
$ws.Range("M16").Value = [double]"21.1309775"
$ws.Range("N16").Value = [double]"42.261955"
$ws.Range("O16").Value = [double]"0.1524111868891667"
$ws.Range("P16").Value = [double]"0.1177170633360173"
$ws.Range("Q16").Value = [double]"2726.940193185593"
$ws.Range("R16").Value = [double]"16361.64115911356"
$ws.Range("S16").Value = [double]"0.04897020159267734"
$ws.Range("T16").Value = [double]"0.03979834083447727"

# Row 17
$ws.Range("E17").Value = [double]"3"
$ws.Range("G17").Value = [double]"158.6435343333333"
$ws.Range("H17").Value = [double]"475.930603"
$ws.Range("I17").Value = [double]"0.3949857228129294"
$ws.Range("J17").Value = [double]"0.4156156560277983"
$ws.Range("K17").Value = [double]"2"
$ws.Range("M17").Value = [double]"35.789624"
$ws.Range("N17").Value = [double]"71.57924800000001"
$ws.Range("O17").Value = [double]"0.258139457682779"
$ws.Range("P17").Value = [double]"0.1993778771086309"
$ws.Range("Q17").Value = [double]"5677.792443821091"
$ws.Range("R17").Value = [double]"34066.75466292654"
$ws.Range("S17").Value = [double]"0.10196140027937"
$ws.Range("T17").Value = [double]"0.08286456719193339"

# Row 18
$ws.Range("E18").Value = [double]"3"
$ws.Range("G18").Value = [double]"158.6435343333333"
$ws.Range("H18").Value = [double]"475.930603"
$ws.Range("I18").Value = [double]"0.3949857228129294"
$ws.Range("J18").Value = [double]"0.4156156560277983"
$ws.Range("K18").Value = [double]"3"
$ws.Range("M18").Value = [double]"60.113367"
$ws.Range("N18").Value = [double]"180.340101"
$ws.Range("O18").Value = [double]"0.4335790718803266"
$ws.Range("P18").Value = [double]"0.5023219368682956"
$ws.Range("Q18").Value = [double]"9536.597001556767"
$ws.Range("R18").Value = [double]"85829.3730140109"
$ws.Range("S18").Value = [double]"0.1712575431032099"
$ws.Range("T18").Value = [double]"0.208772861328671"

# Row 19
$ws.Range("E19").Value = [double]"3"
$ws.Range("G19").Value = [double]"158.6435343333333"
$ws.Range("H19").Value = [double]"475.930603"
$ws.Range("I19").Value = [double]"0.3949857228129294"
$ws.Range("J19").Value = [double]"0.4156156560277983"
$ws.Range("K19").Value = [double]"3"
$ws.Range("M19").Value = [double]"0.2072186666666667"
$ws.Range("N19").Value = [double]"0.621656"
$ws.Range("O19").Value = [double]"0.001494603973349423"
$ws.Range("P19").Value = [double]"0.001731569652308208"
$ws.Range("Q19").Value = [double]"32.87390165984088"
$ws.Range("R19").Value = [double]"295.865114938568"
$ws.Range("S19").Value = [double]"0.0005903472307324982"
$ws.Range("T19").Value = [double]"0.0007196674570019027"

# Row 20
$ws.Range("E20").Value = [double]"3"
$ws.Range("G20").Value = [double]"158.6435343333333"
$ws.Range("H20").Value = [double]"475.930603"
$ws.Range("I20").Value = [double]"0.3949857228129294"
$ws.Range("J20").Value = [double]"0.4156156560277983"
$ws.Range("K20").Value = [double]"3"
$ws.Range("M20").Value = [double]"21.40334366666667"
$ws.Range("N20").Value = [double]"64.210031"
$ws.Range("O20").Value = [double]"0.1543756795743782"
$ws.Range("P20").Value = [double]"0.178851553034748"
$ws.Range("Q20").Value = [double]"3395.502085830966"
$ws.Range("R20").Value = [double]"30559.51877247869"
$ws.Range("S20").Value = [double]"0.06097618938142295"
$ws.Range("T20").Value = [double]"0.07433350554612735"

# Row 21
$ws.Range("E21").Value = [double]"3"
$ws.Range("G21").Value = [double]"158.6435343333333"
$ws.Range("H21").Value = [double]"475.930603"
$ws.Range("I21").Value = [double]"0.3949857228129294"
$ws.Range("J21").Value = [double]"0.4156156560277983"
$ws.Range("K21").Value = [double]"2"
$ws.Range("M21").Value = [double]"21.1309775"
$ws.Range("N21").Value = [double]"42.261955"
$ws.Range("O21").Value = [double]"0.1524111868891667"
$ws.Range("P21").Value = [double]"0.1177170633360173"
$ws.Range("Q21").Value = [double]"3352.292954518144"
$ws.Range("R21").Value = [double]"20113.75772710886"
$ws.Range("S21").Value = [double]"0.06020024281819396"
$ws.Range("T21").Value = [double]"0.04892505450406471"

# Row 22
$ws.Range("E22").Value = [double]"3"
$ws.Range("G22").Value = [double]"41.810285"
$ws.Range("H22").Value = [double]"125.430855"
$ws.Range("I22").Value = [double]"0.1040979432987182"
$ws.Range("J22").Value = [double]"0.1095349337872956"
$ws.Range("K22").Value = [double]"2"
$ws.Range("M22").Value = [double]"35.789624"
$ws.Range("N22").Value = [double]"71.57924800000001"
$ws.Range("O22").Value = [double]"0.258139457682779"
$ws.Range("P22").Value = [double]"0.1993778771086309"
$ws.Range("Q22").Value = [double]"1496.37437948284"
$ws.Range("R22").Value = [double]"8978.246276897042"
$ws.Range("S22").Value = [double]"0.02687178662902378"
$ws.Range("T22").Value = [double]"0.02183884256774544"

# Row 23
$ws.Range("E23").Value = [double]"3"
$ws.Range("G23").Value = [double]"41.810285"
$ws.Range("H23").Value = [double]"125.430855"
$ws.Range("I23").Value = [double]"0.1040979432987182"
$ws.Range("J23").Value = [double]"0.1095349337872956"
$ws.Range("K23").Value = [double]"3"
$ws.Range("M23").Value = [double]"60.113367"
$ws.Range("N23").Value = [double]"180.340101"
$ws.Range("O23").Value = [double]"0.4335790718803266"
$ws.Range("P23").Value = [double]"0.5023219368682956"
$ws.Range("Q23").Value = [double]"2513.357006579595"
$ws.Range("R23").Value = [double]"22620.21305921636"
$ws.Range("S23").Value = [double]"0.04513468964010909"
$ws.Range("T23").Value = [double]"0.05502180009477483"

# Row 24
$ws.Range("E24").Value = [double]"3"
$ws.Range("G24").Value = [double]"41.810285"
$ws.Range("H24").Value = [double]"125.430855"
$ws.Range("I24").Value = [double]"0.1040979432987182"
$ws.Range("J24").Value = [double]"0.1095349337872956"
$ws.Range("K24").Value = [double]"3"
$ws.Range("M24").Value = [double]"0.2072186666666667"
$ws.Range("N24").Value = [double]"0.621656"
$ws.Range("O24").Value = [double]"0.001494603973349423"
$ws.Range("P24").Value = [double]"0.001731569652308208"
$ws.Range("Q24").Value = [double]"8.663871510653333"
$ws.Range("R24").Value = [double]"77.97484359588"
$ws.Range("S24").Value = [double]"0.0001555851996717671"
$ws.Range("T24").Value = [double]"0.00018966736721367"

# Row 25
$ws.Range("E25").Value = [double]"3"
$ws.Range("G25").Value = [double]"41.810285"
$ws.Range("H25").Value = [double]"125.430855"
$ws.Range("I25").Value = [double]"0.1040979432987182"
$ws.Range("J25").Value = [double]"0.1095349337872956"
$ws.Range("K25").Value = [double]"3"
$ws.Range("M25").Value = [double]"21.40334366666667"
$ws.Range("N25").Value = [double]"64.210031"
$ws.Range("O25").Value = [double]"0.1543756795743782"
$ws.Range("P25").Value = [double]"0.178851553034748"
$ws.Range("Q25").Value = [double]"894.8798986562784"
$ws.Range("R25").Value = [double]"8053.919087906505"
$ws.Range("S25").Value = [double]"0.01607019073903471"
$ws.Range("T25").Value = [double]"0.0195904930194161"

# Row 26
$ws.Range("E26").Value = [double]"3"
$ws.Range("G26").Value = [double]"41.810285"
$ws.Range("H26").Value = [double]"125.430855"
$ws.Range("I26").Value = [double]"0.1040979432987182"
$ws.Range("J26").Value = [double]"0.1095349337872956"
$ws.Range("K26").Value = [double]"2"
$ws.Range("M26").Value = [double]"21.1309775"
$ws.Range("N26").Value = [double]"42.261955"
$ws.Range("O26").Value = [double]"0.1524111868891667"
$ws.Range("P26").Value = [double]"0.1177170633360173"
$ws.Range("Q26").Value = [double]"883.4921916035876"
$ws.Range("R26").Value = [double]"5300.953149621525"
$ws.Range("S26").Value = [double]"0.01586569109087881"
$ws.Range("T26").Value = [double]"0.01289413073814553"

# Row 27
$ws.Range("E27").Value = [double]"2"
$ws.Range("G27").Value = [double]"12.780919"
$ws.Range("H27").Value = [double]"25.561838"
$ws.Range("I27").Value = [double]"0.03182153341857176"
$ws.Range("J27").Value = [double]"0.02232237221704002"
$ws.Range("K27").Value = [double]"2"
$ws.Range("M27").Value = [double]"35.789624"
$ws.Range("N27").Value = [double]"71.57924800000001"
$ws.Range("O27").Value = [double]"0.258139457682779"
$ws.Range("P27").Value = [double]"0.1993778771086309"
$ws.Range("Q27").Value = [double]"457.4242853844561"
$ws.Range("R27").Value = [double]"1829.697141537824"
$ws.Range("S27").Value = [double]"0.008214393379304541"
$ws.Range("T27").Value = [double]"0.004450587184662122"

# Row 28
$ws.Range("E28").Value = [double]"2"
$ws.Range("G28").Value = [double]"12.780919"
$ws.Range("H28").Value = [double]"25.561838"
$ws.Range("I28").Value = [double]"0.03182153341857176"
$ws.Range("J28").Value = [double]"0.02232237221704002"
$ws.Range("K28").Value = [double]"3"
$ws.Range("M28").Value = [double]"60.113367"
$ws.Range("N28").Value = [double]"180.340101"
$ws.Range("O28").Value = [double]"0.4335790718803266"
$ws.Range("P28").Value = [double]"0.5023219368682956"
$ws.Range("Q28").Value = [double]"768.3040744442731"
$ws.Range("R28").Value = [double]"4609.824446665639"
$ws.Range("S28").Value = [double]"0.01379715092543314"
$ws.Range("T28").Value = [double]"0.01121301724755857"

# Row 29
$ws.Range("E29").Value = [double]"2"
$ws.Range("G29").Value = [double]"12.780919"
$ws.Range("H29").Value = [double]"25.561838"
$ws.Range("I29").Value = [double]"0.03182153341857176"
$ws.Range("J29").Value = [double]"0.02232237221704002"
$ws.Range("K29").Value = [double]"3"
$ws.Range("M29").Value = [double]"0.2072186666666667"
$ws.Range("N29").Value = [double]"0.621656"
$ws.Range("O29").Value = [double]"0.001494603973349423"
$ws.Range("P29").Value = [double]"0.001731569652308208"
$ws.Range("Q29").Value = [double]"2.648444993954667"
$ws.Range("R29").Value = [double]"15.890669963728"
$ws.Range("S29").Value = [double]"4.75605902854688E-05"
$ws.Range("T29").Value = [double]"3.865274229855441E-05"

# Row 30
$ws.Range("E30").Value = [double]"2"
$ws.Range("G30").Value = [double]"12.780919"
$ws.Range("H30").Value = [double]"25.561838"
$ws.Range("I30").Value = [double]"0.03182153341857176"
$ws.Range("J30").Value = [double]"0.02232237221704002"
$ws.Range("K30").Value = [double]"3"
$ws.Range("M30").Value = [double]"21.40334366666667"
$ws.Range("N30").Value = [double]"64.210031"
$ws.Range("O30").Value = [double]"0.1543756795743782"
$ws.Range("P30").Value = [double]"0.178851553034748"
$ws.Range("Q30").Value = [double]"273.5544017328297"
$ws.Range("R30").Value = [double]"1641.326410396978"
$ws.Range("S30").Value = [double]"0.004912470846590801"
$ws.Range("T30").Value = [double]"0.003992390938437318"

# Row 31
$ws.Range("E31").Value = [double]"2"
$ws.Range("G31").Value = [double]"12.780919"
$ws.Range("H31").Value = [double]"25.561838"
$ws.Range("I31").Value = [double]"0.03182153341857176"
$ws.Range("J31").Value = [double]"0.02232237221704002"
$ws.Range("K31").Value = [double]"2"
$ws.Range("M31").Value = [double]"21.1309775"
$ws.Range("N31").Value = [double]"42.261955"
$ws.Range("O31").Value = [double]"0.1524111868891667"
$ws.Range("P31").Value = [double]"0.1177170633360173"
$ws.Range("Q31").Value = [double]"270.0733118183225"
$ws.Range("R31").Value = [double]"1080.29324727329"
$ws.Range("S31").Value = [double]"0.004849957676957803"
$ws.Range("T31").Value = [double]"0.002627724104083453"
